$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Signature" column (E) was merged E6:E16 and carried the text
# "signature" in E6. Split it into individual bordered (unmerged) cells
# for rows 6-15, matching the blank bordered cells already used in
# columns A-D for those rows, and drop the signature text. Row 16's E
# cell is removed entirely (row 16 keeps only its bottom-border style
# in columns A-D).

# 1. Break the E6:E16 merge so each row in column E becomes an
#    independently addressable cell.
$ws.Range("E6:E16").UnMerge()

# 2. Clear the "signature" text out of E6.
$ws.Range("E6").ClearContents()

# 3. Give E6:E15 the same bordered (no-wrap) formatting already used by
#    the blank filler cells in columns A-D (e.g. D7), so the new column E
#    cells visually match the rest of the table rows.
$ws.Range("D7").Copy()
$ws.Range("E6:E15").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 4. Row 16 no longer has an E cell at all - remove it completely
#    (contents + formatting), unlike rows 6-15 which keep a blank,
#    bordered placeholder cell.
$ws.Range("E16").Clear()
